$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeekApr8")
$ws.Activate()

$xlPasteFormats = -4122

# --- Row 3: "Work on test plan" -- add completion date + actual effort ---
$ws.Range("F3").Value = 43203
$ws.Range("G3").Value = 2

# --- Row 4: "Update buisness model..." -- add completion date + actual effort ---
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F4").Value = 43203
$ws.Range("G4").Value = 1

# --- Row 6: "Create design alternatives" -- add percent complete/completion date/actual effort ---
$ws.Range("E6").Value = 100
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F6").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F6").Value = 43202
$ws.Range("G6").Value = 1

# --- Row 7: "Create design alternitive document" -- same kind of update ---
$ws.Range("E7").Value = 100
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F7").Value = 43202
$ws.Range("G7").Value = 2

# --- Row 9: "Create PlayGame.java" -- effort increased, progress recorded ---
$ws.Range("C9").Value = 1.5
$ws.Range("E9").Value = 15
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F9").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G9").Value = 2

# --- Row 10: "Create Hand.java" -- effort increased, progress recorded ---
$ws.Range("C10").Value = 1.5
$ws.Range("E10").Value = 95
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F10").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G10").Value = 1

# --- Rows 11-13 ("Create BuildPhase.java" / "Create SpacePhase.java" / "Create Scorecard.java") removed from plan ---
$ws.Range("A11:G13").ClearContents() | Out-Null

# --- Selection moved to E9 ---
$ws.Range("E9").Select() | Out-Null

# --- Page orientation set to portrait ---
$ws.PageSetup.Orientation = 1
